$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.612.94'
$ws.Range("E2").Value = '  -0.85%  '
$ws.Range("D3").Value = '3.513.92'
$ws.Range("E3").Value = '  -1.47%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'586.82"
$ws.Range("E5").Value = '  -2.19%  '
$ws.Range("D6").Value = "'133.09"
$ws.Range("E6").Value = '  -0.38%  '
$ws.Range("D7").Value = '3.513.46'
$ws.Range("E7").Value = '  -1.47%  '
$ws.Range("E9").Value = '  -0.82%  '
$ws.Range("E10").Value = '  +0.96%  '
$ws.Range("D11").Value = "'7.15"
$ws.Range("E11").Value = '  +0.20%  '
$ws.Range("E12").Value = '  +0.01%  '
$ws.Range("D13").Value = '4.113.59'
$ws.Range("E13").Value = '  -1.45%  '
$ws.Range("D14").Value = "'27.93"
$ws.Range("E14").Value = '  +3.22%  '
$ws.Range("E15").Value = '  -1.45%  '
$ws.Range("E16").Value = '  +0.63%  '
$ws.Range("D17").Value = '3.513.69'
$ws.Range("E17").Value = '  -1.49%  '
$ws.Range("D18").Value = '64.609.94'
$ws.Range("D19").Value = "'10.02"
$ws.Range("E19").Value = '  +0.71%  '
$ws.Range("D20").Value = "'14.26"
$ws.Range("E20").Value = '  -1.42%  '
$ws.Range("D21").Value = "'5.72"
$ws.Range("E21").Value = '  -2.37%  '
$ws.Range("D22").Value = "'393.34"
$ws.Range("D23").Value = "'0.579"
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("E24").Value = '  -1.46%  '
$ws.Range("D25").Value = "'74.32"
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("E27").Value = '  -2.73%  '
$ws.Range("E28").Value = '  -0.90%  '
$ws.Range("D29").Value = "'7.52"
$ws.Range("E29").Value = '  -4.13%  '
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("E31").Value = '  -1.21%  '
$ws.Range("D32").Value = "'8.24"
$ws.Range("E32").Value = '  -4.55%  '
$ws.Range("D33").Value = '3.517.54'
$ws.Range("E33").Value = '  -1.43%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("E35").Value = '  -0.29%  '
$ws.Range("E36").Value = '  -0.70%  '
$ws.Range("E37").Value = '  +5.36%  '
$ws.Range("E38").Value = '  +3.93%  '
$ws.Range("D39").Value = "'171.32"
$ws.Range("E39").Value = '  +0.36%  '
$ws.Range("D40").Value = "'6.98"
$ws.Range("D41").Value = "'0.0813"
$ws.Range("E41").Value = '  -0.47%  '
$ws.Range("D42").Value = "'26.74"
$ws.Range("E42").Value = '  +0.49%  '
$ws.Range("E43").Value = '  -1.46%  '
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("D45").Value = "'42.36"
$ws.Range("E45").Value = '  -1.62%  '
$ws.Range("D46").Value = "'1.21"
$ws.Range("E46").Value = '  -3.54%  '
$ws.Range("E47").Value = '  -0.87%  '
$ws.Range("E48").Value = '  -0.28%  '
$ws.Range("D49").Value = '2.475.94'
$ws.Range("E49").Value = '  +1.48%  '
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("D51").Value = "'0.914"
$ws.Range("E51").Value = '  +5.26%  '
